# Split Hoja2 ("Numero de facturas") into two sheets:
#  - the first 16 invoice numbers move to a brand-new "Hoja3" sheet
#  - the remaining 18 rows shift up to become the new Hoja2 (rows 1-18)

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Hoja2")

# Add the new sheet after the last existing sheet, so it becomes Hoja3
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Hoja3"

# Move the first 16 rows (A1:A16) of Hoja2 into the new Hoja3
$ws2.Range("A1:A16").Copy() | Out-Null
$ws3.Range("A1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws2.Range("A1:A16").EntireRow.Delete() | Out-Null

# Restore the selection/view state seen on each sheet
$ws3.Range("A1:B16").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D8").Select() | Out-Null
